$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row: add F1 = REX_DEF ---
$ws.Range("F1").Value = 'REX_DEF'
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# --- Rewrite existing rows 2-7 and add new rows 8-17 with final data ---
# (Existing rows 2-7 get their B/C/D/E values replaced per the updated mapping,
#  plus a new F column; rows 8-17 are brand new.)

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 'http://purl.obolibrary.org/obo/ENVO_01000723'
$ws.Range("B7").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B2"), 'http://purl.obolibrary.org/obo/ENVO_01000723')
$ws.Range("B7").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("C2").Value = '{''label'': ''melting'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_01000723''}'
$ws.Range("D2").Value = 'http://purl.obolibrary.org/obo/REX_0000177'
$ws.Range("D7").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D2"), 'http://purl.obolibrary.org/obo/REX_0000177')
$ws.Range("D7").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").Value = '{''label'': ''melting''}'
$ws.Range("F2").Value = '[]'

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 'http://purl.obolibrary.org/obo/ENVO_01000841'
$ws.Range("B7").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B3"), 'http://purl.obolibrary.org/obo/ENVO_01000841')
$ws.Range("B7").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C3").Value = '{''label'': ''pyrolysis'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_01000841''}'
$ws.Range("D3").Value = 'http://purl.obolibrary.org/obo/REX_0000404'
$ws.Range("D7").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D3"), 'http://purl.obolibrary.org/obo/REX_0000404')
$ws.Range("D7").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("E3").Value = '{''label'': ''pyrolysis''}'
$ws.Range("F3").Value = '[]'

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 'http://purl.obolibrary.org/obo/ENVO_01000875'
$ws.Range("B7").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B4"), 'http://purl.obolibrary.org/obo/ENVO_01000875')
$ws.Range("B7").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("C4").Value = '{''label'': ''precipitation'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_01000875''}'
$ws.Range("D4").Value = 'http://purl.obolibrary.org/obo/REX_0000182'
$ws.Range("D7").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D4"), 'http://purl.obolibrary.org/obo/REX_0000182')
$ws.Range("D7").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("E4").Value = '{''label'': ''precipitation''}'
$ws.Range("F4").Value = '[]'

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 'http://purl.obolibrary.org/obo/ENVO_02500034'
$ws.Range("B7").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B5"), 'http://purl.obolibrary.org/obo/ENVO_02500034')
$ws.Range("B7").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("C5").Value = '{''label'': ''evaporation'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_02500034''}'
$ws.Range("D5").Value = 'http://purl.obolibrary.org/obo/REX_0000178'
$ws.Range("D7").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D5"), 'http://purl.obolibrary.org/obo/REX_0000178')
$ws.Range("D7").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = '{''label'': ''evaporation''}'
$ws.Range("F5").Value = '[]'

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 'http://purl.obolibrary.org/obo/ENVO_01000727'
$ws.Range("B7").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B6"), 'http://purl.obolibrary.org/obo/ENVO_01000727')
$ws.Range("B7").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("C6").Value = '{''label'': ''sublimation'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_01000727''}'
$ws.Range("D6").Value = 'http://purl.obolibrary.org/obo/REX_0000180'
$ws.Range("D7").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D6"), 'http://purl.obolibrary.org/obo/REX_0000180')
$ws.Range("D7").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").Value = '{''label'': ''sublimation''}'
$ws.Range("F6").Value = '[]'

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 'http://purl.obolibrary.org/obo/ENVO_01000840'
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B7"), 'http://purl.obolibrary.org/obo/ENVO_01000840')
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("C7").Value = '{''label'': ''thermolysis'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_01000840''}'
$ws.Range("D7").Value = 'http://purl.obolibrary.org/obo/REX_0000086'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D7"), 'http://purl.obolibrary.org/obo/REX_0000086')
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E7").Value = '{''label'': ''thermolysis''}'
$ws.Range("F7").Value = '[]'

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B8").Value = 'http://purl.obolibrary.org/obo/ENVO_01000913'
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B8"), 'http://purl.obolibrary.org/obo/ENVO_01000913')
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("C8").Value = '{''label'': ''nucleation'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_01000913''}'
$ws.Range("D8").Value = 'http://purl.obolibrary.org/obo/REX_0000190'
$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D8"), 'http://purl.obolibrary.org/obo/REX_0000190')
$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("E8").Value = '{''label'': ''nucleation''}'
$ws.Range("F8").Value = '[]'

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("A7").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B9").Value = 'http://purl.obolibrary.org/obo/ENVO_01000931'
$ws.Range("B7").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B9"), 'http://purl.obolibrary.org/obo/ENVO_01000931')
$ws.Range("B7").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C9").Value = '{''label'': ''migration'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_01000931''}'
$ws.Range("D9").Value = 'http://purl.obolibrary.org/obo/REX_0000374'
$ws.Range("D7").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D9"), 'http://purl.obolibrary.org/obo/REX_0000374')
$ws.Range("D7").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("E9").Value = '{''label'': ''migration''}'
$ws.Range("F9").Value = '[]'

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("A7").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").Value = 'http://purl.obolibrary.org/obo/ENVO_01001259'
$ws.Range("B7").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B10"), 'http://purl.obolibrary.org/obo/ENVO_01001259')
$ws.Range("B7").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("C10").Value = '{''label'': ''transport'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_01001259''}'
$ws.Range("D10").Value = 'http://purl.obolibrary.org/obo/REX_0000458'
$ws.Range("D7").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D10"), 'http://purl.obolibrary.org/obo/REX_0000458')
$ws.Range("D7").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E10").Value = '{''label'': ''transport''}'
$ws.Range("F10").Value = '[]'

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("A7").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B11").Value = 'http://purl.obolibrary.org/obo/ENVO_01001261'
$ws.Range("B7").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B11"), 'http://purl.obolibrary.org/obo/ENVO_01001261')
$ws.Range("B7").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C11").Value = '{''label'': ''transport'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_01001261''}'
$ws.Range("D11").Value = 'http://purl.obolibrary.org/obo/REX_0000458'
$ws.Range("D7").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D11"), 'http://purl.obolibrary.org/obo/REX_0000458')
$ws.Range("D7").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E11").Value = '{''label'': ''transport''}'
$ws.Range("F11").Value = '[]'

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("A7").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("B12").Value = 'http://purl.obolibrary.org/obo/ENVO_01001262'
$ws.Range("B7").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B12"), 'http://purl.obolibrary.org/obo/ENVO_01001262')
$ws.Range("B7").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("C12").Value = '{''label'': ''transport'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_01001262''}'
$ws.Range("D12").Value = 'http://purl.obolibrary.org/obo/REX_0000458'
$ws.Range("D7").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D12"), 'http://purl.obolibrary.org/obo/REX_0000458')
$ws.Range("D7").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E12").Value = '{''label'': ''transport''}'
$ws.Range("F12").Value = '[]'

# Row 13
$ws.Range("A13").Value = 11
$ws.Range("A7").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("B13").Value = 'http://purl.obolibrary.org/obo/ENVO_03400014'
$ws.Range("B7").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B13"), 'http://purl.obolibrary.org/obo/ENVO_03400014')
$ws.Range("B7").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C13").Value = '{''label'': ''precipitation'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_03400014''}'
$ws.Range("D13").Value = 'http://purl.obolibrary.org/obo/REX_0000182'
$ws.Range("D7").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D13"), 'http://purl.obolibrary.org/obo/REX_0000182')
$ws.Range("D7").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("E13").Value = '{''label'': ''precipitation''}'
$ws.Range("F13").Value = '[]'

# Row 14
$ws.Range("A14").Value = 12
$ws.Range("A7").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("B14").Value = 'http://purl.obolibrary.org/obo/ENVO_03501325'
$ws.Range("B7").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B14"), 'http://purl.obolibrary.org/obo/ENVO_03501325')
$ws.Range("B7").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C14").Value = '{''label'': ''transport'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_03501325''}'
$ws.Range("D14").Value = 'http://purl.obolibrary.org/obo/REX_0000458'
$ws.Range("D7").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D14"), 'http://purl.obolibrary.org/obo/REX_0000458')
$ws.Range("D7").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = '{''label'': ''transport''}'
$ws.Range("F14").Value = '[]'

# Row 15
$ws.Range("A15").Value = 13
$ws.Range("A7").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("B15").Value = 'http://purl.obolibrary.org/obo/ENVO_03501117'
$ws.Range("B7").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B15"), 'http://purl.obolibrary.org/obo/ENVO_03501117')
$ws.Range("B7").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C15").Value = '{''label'': ''transport'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_03501117''}'
$ws.Range("D15").Value = 'http://purl.obolibrary.org/obo/REX_0000458'
$ws.Range("D7").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D15"), 'http://purl.obolibrary.org/obo/REX_0000458')
$ws.Range("D7").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = '{''label'': ''transport''}'
$ws.Range("F15").Value = '[]'

# Row 16
$ws.Range("A16").Value = 14
$ws.Range("A7").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("B16").Value = 'http://purl.obolibrary.org/obo/ENVO_06105021'
$ws.Range("B7").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B16"), 'http://purl.obolibrary.org/obo/ENVO_06105021')
$ws.Range("B7").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C16").Value = '{''label'': ''adsorption'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_06105021''}'
$ws.Range("D16").Value = 'http://purl.obolibrary.org/obo/REX_0000198'
$ws.Range("D7").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D16"), 'http://purl.obolibrary.org/obo/REX_0000198')
$ws.Range("D7").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = '{''label'': ''adsorption''}'
$ws.Range("F16").Value = '[]'

# Row 17
$ws.Range("A17").Value = 15
$ws.Range("A7").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B17").Value = 'http://purl.obolibrary.org/obo/ENVO_09000028'
$ws.Range("B7").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B17"), 'http://purl.obolibrary.org/obo/ENVO_09000028')
$ws.Range("B7").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C17").Value = '{''label'': ''evaporation'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_09000028''}'
$ws.Range("D17").Value = 'http://purl.obolibrary.org/obo/REX_0000178'
$ws.Range("D7").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D17"), 'http://purl.obolibrary.org/obo/REX_0000178')
$ws.Range("D7").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = '{''label'': ''evaporation''}'
$ws.Range("F17").Value = '[]'

